# Update the "Well_1" label (the second occurrence, shape "TextBox 154")
# in the file-structure diagram to "Well_n", matching the author's
# "Updating expected file structure image" commit.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item("TextBox 154")
$shape.TextFrame.TextRange.Text = "Well_n"
